$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.25
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("I2").Value = 0.3333333333333333
$ws.Range("U2").Value = 1
$ws.Range("W2").Value = 1
$ws.Range("C3").Value = 0.3333333333333333
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("I3").Value = 0.3333333333333333
$ws.Range("U3").Value = 0.3333333333333333
$ws.Range("W3").Value = 0.6666666666666666
$ws.Range("C4").Value = 0.2857142857142858
$ws.Range("F4").Value = 0.4444444444444444
$ws.Range("I4").Value = 0.3333333333333333
$ws.Range("U4").Value = 0.5
$ws.Range("W4").Value = 0.8
$ws.Range("C5").Value = 0.3125
$ws.Range("F5").Value = 0.5555555555555555
$ws.Range("I5").Value = 0.3333333333333333
$ws.Range("U5").Value = 0.3846153846153846
$ws.Range("W5").Value = 0.7142857142857142
$ws.Range("C6").Value = 0.4702019977678391
$ws.Range("F6").Value = 0.5113881456198478
$ws.Range("I6").Value = 0.7452525342261976
$ws.Range("U6").Value = 0.7452525342261976
$ws.Range("W6").Value = 0.5766666455144387
$ws.Range("I7").Value = $true
$ws.Range("U7").Value = $true
$ws.Range("C8").Value = $true
$ws.Range("F8").Value = $true
$ws.Range("I8").Value = $true
$ws.Range("U8").Value = $true
$ws.Range("W8").Value = $true
$ws.Range("C9").Value = $true
$ws.Range("F9").Value = $true
$ws.Range("I9").Value = $true
$ws.Range("U9").Value = $true
$ws.Range("W9").Value = $true
$ws.Range("C10").Value = 2
$ws.Range("F10").Value = 2
$ws.Range("I10").Value = 1
$ws.Range("U10").Value = 1
$ws.Range("W10").Value = 2
